# fix: remove "Notice u/s 94 BNSS" heading paragraph from this
# "Put on Hold" bank letter template (the money-release bullet-spacing
# fix from the same commit lives in a different template file and does
# not apply to this document).

$d = $word.ActiveDocument

# Walk paragraphs back-to-front (safe if more than one ever matches) and
# delete the whole paragraph (text + its paragraph mark) whose text is
# the "Notice u/s 94 BNSS" heading.
for ($i = $d.Paragraphs.Count; $i -ge 1; $i--) {
    $para = $d.Paragraphs.Item($i)
    if ($para.Range.Text -like "*Notice u/s 94 BNSS*") {
        $para.Range.Delete()
    }
}
